$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap match data (columns F:V) between paired rows; columns A:E (index/meta) stay put.
$rA = $ws.Range("F8:V8")
$rB = $ws.Range("F9:V9")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

$rA = $ws.Range("F18:V18")
$rB = $ws.Range("F19:V19")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

$rA = $ws.Range("F29:V29")
$rB = $ws.Range("F30:V30")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

$rA = $ws.Range("F31:V31")
$rB = $ws.Range("F32:V32")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

$rA = $ws.Range("F49:V49")
$rB = $ws.Range("F50:V50")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

$rA = $ws.Range("F51:V51")
$rB = $ws.Range("F52:V52")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

$rA = $ws.Range("F63:V63")
$rB = $ws.Range("F64:V64")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

$rA = $ws.Range("F68:V68")
$rB = $ws.Range("F69:V69")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

$rA = $ws.Range("F75:V75")
$rB = $ws.Range("F76:V76")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

$rA = $ws.Range("F79:V79")
$rB = $ws.Range("F80:V80")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

$rA = $ws.Range("F91:V91")
$rB = $ws.Range("F92:V92")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

$rA = $ws.Range("F115:V115")
$rB = $ws.Range("F116:V116")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value2 = $vB
$rB.Value2 = $vA

# Append two new match rows (136, 137), copying row 135 formatting first.
$ws.Range("A135:V135").Copy()
$ws.Range("A136:V137").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(136,1).Value2 = 135
$ws.Cells.Item(136,2).Value2 = "turkey"
$ws.Cells.Item(136,3).Value2 = "super-lig"
$ws.Cells.Item(136,4).Value2 = "2023-2024"
$ws.Cells.Item(136,5).Value2 = 45262.70833333334
$ws.Cells.Item(136,6).Value2 = "Adana Demirspor"
$ws.Cells.Item(136,7).Value2 = 2
$ws.Cells.Item(136,8).Value2 = "Samsunspor"
$ws.Cells.Item(136,9).Value2 = 3
$ws.Cells.Item(136,10).Value2 = 1.85
$ws.Cells.Item(136,11).Value2 = "26/11/2023 14:13"
$ws.Cells.Item(136,12).Value2 = 1.82
$ws.Cells.Item(136,13).Value2 = "02/12/2023 16:57"
$ws.Cells.Item(136,14).Value2 = 4.15
$ws.Cells.Item(136,15).Value2 = "26/11/2023 14:13"
$ws.Cells.Item(136,16).Value2 = 3.9
$ws.Cells.Item(136,17).Value2 = "02/12/2023 16:57"
$ws.Cells.Item(136,18).Value2 = 3.84
$ws.Cells.Item(136,19).Value2 = "26/11/2023 14:13"
$ws.Cells.Item(136,20).Value2 = 4.49
$ws.Cells.Item(136,21).Value2 = "02/12/2023 16:57"
$ws.Cells.Item(136,22).Value2 = "https://www.betexplorer.com/football/turkey/super-lig/adanademirspor-samsunspor/txLELtQo/"

$ws.Cells.Item(137,1).Value2 = 136
$ws.Cells.Item(137,2).Value2 = "turkey"
$ws.Cells.Item(137,3).Value2 = "super-lig"
$ws.Cells.Item(137,4).Value2 = "2023-2024"
$ws.Cells.Item(137,5).Value2 = 45262.70833333334
$ws.Cells.Item(137,6).Value2 = "Pendikspor"
$ws.Cells.Item(137,7).Value2 = 0
$ws.Cells.Item(137,8).Value2 = "Galatasaray"
$ws.Cells.Item(137,9).Value2 = 2
$ws.Cells.Item(137,10).Value2 = 6.33
$ws.Cells.Item(137,11).Value2 = "27/11/2023 18:12"
$ws.Cells.Item(137,12).Value2 = 9.94
$ws.Cells.Item(137,13).Value2 = "02/12/2023 16:54"
$ws.Cells.Item(137,14).Value2 = 5.29
$ws.Cells.Item(137,15).Value2 = "27/11/2023 18:12"
$ws.Cells.Item(137,16).Value2 = 6.31
$ws.Cells.Item(137,17).Value2 = "02/12/2023 16:54"
$ws.Cells.Item(137,18).Value2 = 1.43
$ws.Cells.Item(137,19).Value2 = "27/11/2023 18:12"
$ws.Cells.Item(137,20).Value2 = 1.28
$ws.Cells.Item(137,21).Value2 = "02/12/2023 16:54"
$ws.Cells.Item(137,22).Value2 = "https://www.betexplorer.com/football/turkey/super-lig/pendikspor-galatasaray/8lj02JPN/"


Write-Host "Edit applied successfully"
